{"js": "const pairs = [\n  [\"2025-07-28 Monday\", \"2025-07-29 Tuesday\"],\n  [\"454\u00f72=227, 0\", \"300\u00f73=100, 0\"],\n  [\"416\u00f73=138, 2\", \"504\u00f74=126, 0\"],\n  [\"560\u00f73=186, 2\", \"458\u00f75=91, 3\"],\n  [\"971\u00f73=323, 2\", \"501\u00f77=71, 4\"],\n  [\"167\u00f72=83, 1\", \"694\u00f73=231, 1\"],\n  [\"747\u00f73=249, 0\", \"731\u00f72=365, 1\"],\n  [\"371\u00f74=92, 3\", \"409\u00f78=51, 1\"],\n  [\"520\u00f77=74, 2\", \"590\u00f79=65, 5\"],\n  [\"603\u00f79=67, 0\", \"695\u00f72=347, 1\"],\n  [\"558\u00f78=69, 6\", \"205\u00f78=25, 5\"],\n  [\"152\u00f77=21, 5\", \"213\u00f78=26, 5\"],\n  [\"878\u00f76=146, 2\", \"130\u00f77=18, 4\"],\n  [\"607\u00f74=151, 3\", \"180\u00f75=36, 0\"],\n  [\"985\u00f78=123, 1\", \"755\u00f76=125, 5\"],\n  [\"286\u00f74=71, 2\", \"763\u00f79=84, 7\"],\n  [\"733\u00f73=244, 1\", \"129\u00f72=64, 1\"],\n  [\"751\u00f72=375, 1\", \"265\u00f75=53, 0\"],\n  [\"791\u00f78=98, 7\", \"932\u00f73=310, 2\"],\n  [\"702\u00f76=117, 0\", \"151\u00f77=21, 4\"],\n  [\"908\u00f79=100, 8\", \"256\u00f72=128, 0\"],\n  [\"480\u00f75=96, 0\", \"518\u00f78=64, 6\"],\n  [\"654\u00f79=72, 6\", \"835\u00f78=104, 3\"],\n  [\"309\u00f74=77, 1\", \"385\u00f73=128, 1\"],\n  [\"174\u00f72=87, 0\", \"782\u00f77=111, 5\"],\n  [\"598\u00f72=299, 0\", \"870\u00f73=290, 0\"],\n];\n\nconst body = context.document.body;\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  for (const r of results.items) {\n    r.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n    ,@(\"2025-07-28 Monday\", \"2025-07-29 Tuesday\")\n    ,@(\"454\u00f72=227, 0\", \"300\u00f73=100, 0\")\n    ,@(\"416\u00f73=138, 2\", \"504\u00f74=126, 0\")\n    ,@(\"560\u00f73=186, 2\", \"458\u00f75=91, 3\")\n    ,@(\"971\u00f73=323, 2\", \"501\u00f77=71, 4\")\n    ,@(\"167\u00f72=83, 1\", \"694\u00f73=231, 1\")\n    ,@(\"747\u00f73=249, 0\", \"731\u00f72=365, 1\")\n    ,@(\"371\u00f74=92, 3\", \"409\u00f78=51, 1\")\n    ,@(\"520\u00f77=74, 2\", \"590\u00f79=65, 5\")\n    ,@(\"603\u00f79=67, 0\", \"695\u00f72=347, 1\")\n    ,@(\"558\u00f78=69, 6\", \"205\u00f78=25, 5\")\n    ,@(\"152\u00f77=21, 5\", \"213\u00f78=26, 5\")\n    ,@(\"878\u00f76=146, 2\", \"130\u00f77=18, 4\")\n    ,@(\"607\u00f74=151, 3\", \"180\u00f75=36, 0\")\n    ,@(\"985\u00f78=123, 1\", \"755\u00f76=125, 5\")\n    ,@(\"286\u00f74=71, 2\", \"763\u00f79=84, 7\")\n    ,@(\"733\u00f73=244, 1\", \"129\u00f72=64, 1\")\n    ,@(\"751\u00f72=375, 1\", \"265\u00f75=53, 0\")\n    ,@(\"791\u00f78=98, 7\", \"932\u00f73=310, 2\")\n    ,@(\"702\u00f76=117, 0\", \"151\u00f77=21, 4\")\n    ,@(\"908\u00f79=100, 8\", \"256\u00f72=128, 0\")\n    ,@(\"480\u00f75=96, 0\", \"518\u00f78=64, 6\")\n    ,@(\"654\u00f79=72, 6\", \"835\u00f78=104, 3\")\n    ,@(\"309\u00f74=77, 1\", \"385\u00f73=128, 1\")\n    ,@(\"174\u00f72=87, 0\", \"782\u00f77=111, 5\")\n    ,@(\"598\u00f72=299, 0\", \"870\u00f73=290, 0\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n}"}
